$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Initial release"
$ws.Range("F1").Value = "Deprecation release"
$ws.Range("G1").Value = "Removal Date"
